{"js": "// The document has a hidden \"_GoBack\" bookmark (the mark Word leaves at the\n// last edit location) that currently sits inside the paragraph that reads\n// \"...pain management which includes opioid medications |for people with...\",\n// wrapping nothing (bookmarkStart immediately followed by bookmarkEnd).\n//\n// The authored edit relocates that bookmark so it spans from the very start\n// of the \"The CDC Guidelines were originally phrased...\" paragraph all the\n// way through to the very end of the final \"Dr. and Ms. Hollis\" paragraph\n// (i.e. it now wraps the whole closing block of the letter instead of a\n// single empty point mid-sentence). The visible text is unchanged.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// Locate the anchor paragraphs for the bookmark's new position by matching\n// on their (unique) text rather than a hard-coded index.\nlet startParagraph = null;\nlet endParagraph = null;\n\nfor (const p of paragraphs.items) {\n  p.load(\"text\");\n}\nawait context.sync();\n\nfor (const p of paragraphs.items) {\n  const text = p.text;\n  if (text.indexOf(\"The CDC Guidelines were originally phrased\") !== -1) {\n    startParagraph = p;\n  }\n  if (text.indexOf(\"Dr. and Ms. Hollis\") !== -1) {\n    endParagraph = p;\n  }\n}\n\nif (!startParagraph || !endParagraph) {\n  throw new Error(\"Could not locate anchor paragraphs for the _GoBack bookmark move.\");\n}\n\n// Remove the bookmark from wherever it currently lives.\ncontext.document.deleteBookmark(\"_GoBack\");\n\n// Re-insert it so it spans from the start of the CDC-guidelines paragraph to\n// the end of the closing \"Dr. and Ms. Hollis\" paragraph.\nconst rangeStart = startParagraph.getRange(\"Start\");\nconst rangeEnd = endParagraph.getRange(\"End\");\nconst fullRange = rangeStart.expandTo(rangeEnd);\nfullRange.insertBookmark(\"_GoBack\");\n\nawait context.sync();\n", "ps1": "# The document carries a hidden \"_GoBack\" bookmark (the marker Word drops at\n# the last edit location). It currently sits mid-sentence in the paragraph\n# \"...pain management which includes opioid medications |for people with...\",\n# wrapping nothing (bookmarkStart immediately followed by bookmarkEnd).\n#\n# This edit relocates that bookmark so that it spans from the very start of\n# the \"The CDC Guidelines were originally phrased...\" paragraph through to\n# the very end of the closing \"Dr. and Ms. Hollis\" paragraph, i.e. it now\n# wraps the whole tail block of the letter instead of a single empty point\n# mid-sentence. No visible text changes.\n\n$d = $word.ActiveDocument\n\n# Find the paragraph that anchors the new bookmark start.\n$startFind = $d.Content\n$startFind.Find.Execute(\"The CDC Guidelines were originally phrased\") | Out-Null\n$startParagraph = $startFind.Paragraphs(1)\n\n# Find the paragraph that anchors the new bookmark end (the last paragraph\n# of the letter).\n$endFind = $d.Content\n$endFind.Find.Execute(\"Dr. and Ms. Hollis\") | Out-Null\n$endParagraph = $endFind.Paragraphs(1)\n\n# Build a range spanning from the start of the first paragraph to the end of\n# the last one.\n$newBookmarkRange = $d.Range($startParagraph.Range.Start, $endParagraph.Range.End)\n\n# Re-adding a bookmark with the same name (\"_GoBack\") as an existing one\n# replaces it, which both removes it from its old location and places it at\n# the new one.\n$d.Bookmarks.Add(\"_GoBack\", $newBookmarkRange)\n"}
